$wb = $excel.ActiveWorkbook

# 1. Rename the original sheet to "API-Testing"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "API-Testing"

# 2. Duplicate the sheet (Copy places the new sheet right after the source)
$ws1.Copy()
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "API-Testing-Sheet2-Duplicate"

# 3. On the duplicated sheet, update the "tag=grey" cell (M2) so it gets its
#    own shared-string entry (workflow execution now drives the tag value).
$m2 = $ws2.Range("M2")
$m2.Value = "tag=grey"

$m2Run1 = $m2.Characters(1, 4)
$m2Run1.Font.Name = "Inconsolata"
$m2Run1.Font.Bold = $true
$m2Run1.Font.Underline = $true
$m2Run1.Font.Size = 11
$m2Run1.Font.Color = 0

$m2Run2 = $m2.Characters(5, 4)
$m2Run2.Font.Name = "Inconsolata"
$m2Run2.Font.Bold = $true
$m2Run2.Font.Underline = $true
$m2Run2.Font.Size = 11
$m2Run2.Font.Color = 13391121

# 4. On the duplicated sheet, update the findByTags URL cell (C3) likewise.
$c3 = $ws2.Range("C3")
$c3.Value = "https://live.virtualandemo.com/api/pets/findByTags?tags=[tag]"

$c3Run1 = $c3.Characters(1, 56)
$c3Run1.Font.Name = "Inconsolata"
$c3Run1.Font.Bold = $true
$c3Run1.Font.Underline = $true
$c3Run1.Font.Size = 11
$c3Run1.Font.Color = 13391121

$c3Run2 = $c3.Characters(57, 5)
$c3Run2.Font.Name = "Inconsolata"
$c3Run2.Font.Bold = $true
$c3Run2.Font.Underline = $true
$c3Run2.Font.Size = 11
$c3Run2.Font.Color = 0

# 5. The workflow-execution cells no longer use the hyperlink-ish look;
#    switch them to the plain "theme text" style used elsewhere on the sheet.
$ws2.Range("N2").Font.ThemeColor = 1
$ws2.Range("O2").Font.ThemeColor = 1
$ws2.Range("N3").Font.ThemeColor = 1
$ws2.Range("O3").Font.ThemeColor = 1
$ws2.Range("O4").Font.ThemeColor = 1
